# Updates cryptos list values (Price / Volume(1h) columns) to match the
# latest scrape. All D/E cells in this sheet are stored as text, so we
# force a text number format before writing to avoid Excel coercing
# numeric-looking strings (e.g. '559.33') into floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.047.82'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.17%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.396.11'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.89%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '559.33'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.31%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.74'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.08%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.586'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.24%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.65'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.65%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.32%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.83%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '24.61'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.826.86'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.79%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '59.984.18'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.23%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.22%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.405.22'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.46%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.14'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -2.07%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.44%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '322.75'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.76'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.81%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '64.16'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -3.47%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.25%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.49'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.22%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.998'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.26%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.78%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.80%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.36%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '171.04'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.02%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.12%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +7.14%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.25%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -2.34%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.09%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.10%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.60%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '323.30'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.75%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '38.71'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.26%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '147.26'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +6.22%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -3.35%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '19.85'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.28%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0513'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.17%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.89%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.93%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '11.05'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.06%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.32%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.35%  '
